$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$location = " Morumbi - SP "

$jobs = @(
    'Analista de Sistemas Sênior',
    'Analista de Infraestrutura Bilingue',
    'Analista de Teste',
    'Analista Jurídico Pl II',
    'Analista BI',
    'Scrum Master',
    'Tech Lead - Mobile',
    'Agile Master',
    'Desenvolvedor Mobile - Flutter',
    'Engenheiro de Software Java Especialista',
    'Desenvolvedor Móbile Sênior',
    'Desenvolvedor Backend Sr',
    'Engenheiro de Software Sr (.NET Core e/ou Node.js)',
    'Analista SOA.',
    'UX/UI Designer Pleno',
    'Desenvolvedor Back End NodeJS',
    'Scrum Master',
    'Scrum Master',
    'Desenvolvedor .NET',
    'Desenvolvedor .NET/Oracle',
    'Tech Lead',
    'QA com conhecimento em Jest',
    'Desenvolvedor RPA',
    'Desenvolvedor Full Stack',
    'Desenvolvedor FontEnd React',
    'Desenvolvedor Back End Java',
    'Analista de Processos RPA',
    'Analista de NOC Jr III',
    'Dev Peoplesoft - ERP',
    'Estágio Logistica',
    'Estágio em TI',
    'Especialista em QA',
    'Lider Técnico',
    'Analista Service Desk Jr I',
    'Desenvolvedor BackEnd Java - Springboot',
    'Desenvolvedor FullStack'
)

$startRow = 2
for ($i = 0; $i -lt $jobs.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $jobs[$i]
    $ws.Cells.Item($row, 2).Value = $location
}

$ws.Cells.Item(2, 3).Value = 'Teste'
